$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=129; B="Genf"; C="GE"; D=18012; H=3.57; J=44227},
    @{Row=130; B="Waadt"; C="VD"; D=27848; H=3.46; J=44227},
    @{Row=131; B="Wallis"; C="VS"; D=14240; H=4.12; J=44227},
    @{Row=132; B="Freiburg"; C="FR"; D=9580; H=2.98; J=44227},
    @{Row=133; B="Neuenburg"; C="NE"; D=4907; H=2.78; J=44227},
    @{Row=134; B="Jura"; C="JU"; D=2792; H=3.79; J=44227},
    @{Row=135; B="Bern"; C="BE"; D=24385; H=2.35; J=44227},
    @{Row=136; B="Solothurn"; C="SO"; D=15099; H=5.49; J=44227},
    @{Row=137; B="Basel-Stadt"; C="BS"; D=13754; H=7.02; J=44227},
    @{Row=138; B="Basel-Landsch."; C="BL"; D=11364; H=3.93; J=44227},
    @{Row=139; B="Aargau"; C="AG"; D=20034; H=2.92; J=44227},
    @{Row=140; B="Zürich"; C="ZH"; D=46081; H=2.99; J=44227},
    @{Row=141; B="Schaffhausen"; C="SH"; D=5685; H=6.9; J=44227},
    @{Row=142; B="Thurgau"; C="TG"; D=8868; H=3.17; J=44227},
    @{Row=143; B="Appenzell-A."; C="AR"; D=3542; H=6.39; J=44227},
    @{Row=144; B="Appenzell-I."; C="AI"; D=1092; H=6.77; J=44227},
    @{Row=145; B="St. Gallen"; C="SG"; D=19281; H=3.78; J=44227},
    @{Row=146; B="Glarus"; C="GL"; D=2330; H=5.74; J=44227},
    @{Row=147; B="Schwyz"; C="SZ"; D=5165; H=3.22; J=44227},
    @{Row=148; B="Zug"; C="ZG"; D=7431; H=5.82; J=44227},
    @{Row=149; B="Luzern"; C="LU"; D=18076; H=4.38; J=44227},
    @{Row=150; B="Nidwalden"; C="NW"; D=3218; H=7.47; J=44227},
    @{Row=151; B="Obwalden"; C="OW"; D=2156; H=5.68; J=44227},
    @{Row=152; B="Uri"; C="UR"; D=2184; H=5.95; J=44227},
    @{Row=153; B="Graubünden"; C="GR"; D=7339; H=3.69; J=44227},
    @{Row=154; B="Tessin"; C="TI"; D=19737; H=5.62; J=44227}
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B   # Kanton
    $ws.Cells.Item($r.Row, 3).Value = $r.C   # Abk.
    $ws.Cells.Item($r.Row, 4).Value = $r.D   # Total_Impfungen
    $ws.Cells.Item($r.Row, 8).Value = $r.H   # Impfungen pro 100 Einwohner
    $ws.Cells.Item($r.Row, 10).Value = $r.J  # Stand (date serial)
}

# Apply the date number format to the first new "Stand" cell, then
# propagate it via a format-only paste so every new J cell shares a single
# re-used style record instead of each getting its own duplicate xf.
$ws.Range("J129").NumberFormat = "mm-dd-yy"
$ws.Range("J129").Copy()
$ws.Range("J130:J154").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the view/selection state to match the edited workbook
$ws.Range("B155").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 112
